$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 86.40000000000001
$ws.Range("I4").Value = 90.666664
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 90.666664
$ws.Range("L4").Value = 80
$ws.Range("M4").Value = 23.333336
$ws.Range("N4").Value = -308
$ws.Range("H62").Value = 9789.538
$ws.Range("I62").Value = 9751.272000000001
$ws.Range("K62").Value = 9751.272000000001
$ws.Range("M62").Value = -9127.272000000001
$ws.Range("H65").Value = 9789.538
$ws.Range("I65").Value = 9751.272000000001
$ws.Range("K65").Value = 48756.36
$ws.Range("M65").Value = -45636.36
$ws.Range("H103").Value = 500
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 500
$ws.Range("K103").Value = 0
$ws.Range("M103").Value = 1500
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -2672
$ws.Range("H138").Value = 4227.2144
$ws.Range("I138").Value = 4308.3
$ws.Range("J138").Value = 4182.1665
$ws.Range("K138").Value = 12924.9
$ws.Range("L138").Value = 12546.4995
$ws.Range("M138").Value = -7784.900000000001
$ws.Range("N138").Value = -22826.4995

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 20987
$ws.Range("I28").Value = 20987
$ws.Range("K28").Value = 20987
$ws.Range("M28").Value = -20795
$ws.Range("H61").Value = 16687.25
$ws.Range("I61").Value = 9899.799999999999
$ws.Range("K61").Value = 9899.799999999999
$ws.Range("M61").Value = -9687.799999999999
$ws.Range("H99").Value = 20987
$ws.Range("I99").Value = 20987
$ws.Range("K99").Value = 20987
$ws.Range("M99").Value = -17992
$ws.Range("H110").Value = 1975.0667
$ws.Range("I110").Value = 1636.4166
$ws.Range("K110").Value = 1636.4166
$ws.Range("M110").Value = 408.5834
$ws.Range("H136").Value = 16687.25
$ws.Range("I136").Value = 9899.799999999999
$ws.Range("K136").Value = 29699.4
$ws.Range("M136").Value = -27149.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 11411.875
$ws.Range("I54").Value = 11411.875
$ws.Range("K54").Value = 11411.875
$ws.Range("M54").Value = -10927.875
$ws.Range("H82").Value = 49047.707
$ws.Range("I82").Value = 11539.223
$ws.Range("K82").Value = 11539.223
$ws.Range("M82").Value = -11156.223
$ws.Range("H85").Value = 49047.707
$ws.Range("I85").Value = 11539.223
$ws.Range("K85").Value = 11539.223
$ws.Range("M85").Value = -10213.223
$ws.Range("H94").Value = 1766.9412
$ws.Range("I94").Value = 1564.1428
$ws.Range("J94").Value = 2713.3333
$ws.Range("K94").Value = 1564.1428
$ws.Range("L94").Value = 2713.3333
$ws.Range("M94").Value = -1113.1428
$ws.Range("N94").Value = -3615.3333
$ws.Range("H99").Value = 8177.4
$ws.Range("I99").Value = 9221.75
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 9221.75
$ws.Range("L99").Value = 4000
$ws.Range("M99").Value = -7723.75
$ws.Range("N99").Value = -6996
$ws.Range("H105").Value = 4116.3438
$ws.Range("I105").Value = 3610.84
$ws.Range("K105").Value = 3610.84
$ws.Range("M105").Value = -1863.84
$ws.Range("H130").Value = 78332.164
$ws.Range("J130").Value = 78332.164
$ws.Range("L130").Value = 78332.164
$ws.Range("N130").Value = -88372.164
$ws.Range("H132").Value = 84953.914
$ws.Range("J132").Value = 84953.914
$ws.Range("L132").Value = 84953.914
$ws.Range("N132").Value = -95073.914
$ws.Range("H134").Value = 24191.268
$ws.Range("I134").Value = 17255.25
$ws.Range("J134").Value = 26713.455
$ws.Range("K134").Value = 51765.75
$ws.Range("L134").Value = 80140.36500000001
$ws.Range("M134").Value = -49230.75
$ws.Range("N134").Value = -85210.36500000001
$ws.Range("H137").Value = 124999
$ws.Range("J137").Value = 124999
$ws.Range("L137").Value = 124999
$ws.Range("N137").Value = -135199

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 68058.44
$ws.Range("I62").Value = 501349.5
$ws.Range("J62").Value = 6159.7144
$ws.Range("K62").Value = 501349.5
$ws.Range("L62").Value = 6159.7144
$ws.Range("M62").Value = -500725.5
$ws.Range("N62").Value = -7407.7144
$ws.Range("H65").Value = 68058.44
$ws.Range("I65").Value = 501349.5
$ws.Range("J65").Value = 6159.7144
$ws.Range("K65").Value = 2506747.5
$ws.Range("L65").Value = 30798.572
$ws.Range("M65").Value = -2503627.5
$ws.Range("N65").Value = -37038.572
$ws.Range("H86").Value = 3177.92
$ws.Range("I86").Value = 2609.9285
$ws.Range("J86").Value = 3900.818
$ws.Range("K86").Value = 2609.9285
$ws.Range("L86").Value = 3900.818
$ws.Range("M86").Value = -1486.9285
$ws.Range("N86").Value = -6146.818
$ws.Range("H89").Value = 3177.92
$ws.Range("I89").Value = 2609.9285
$ws.Range("J89").Value = 3900.818
$ws.Range("K89").Value = 13049.6425
$ws.Range("L89").Value = 19504.09
$ws.Range("M89").Value = -7433.6425
$ws.Range("N89").Value = -30736.09
$ws.Range("H95").Value = 15860.5
$ws.Range("J95").Value = 15860.5
$ws.Range("L95").Value = 15860.5
$ws.Range("N95").Value = -21352.5
$ws.Range("H122").Value = 10014
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H140").Value = 63799
$ws.Range("J140").Value = 63799
$ws.Range("L140").Value = 63799
$ws.Range("N140").Value = -74159

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 296.33334
$ws.Range("J40").Value = 700
$ws.Range("L40").Value = 2800
$ws.Range("N40").Value = -2938

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H80").Value = 1198
$ws.Range("I80").Value = 1498.5
$ws.Range("J80").Value = 997.6667
$ws.Range("K80").Value = 1498.5
$ws.Range("L80").Value = 997.6667
$ws.Range("M80").Value = -500.5
$ws.Range("N80").Value = -2993.6667
$ws.Range("H83").Value = 1198
$ws.Range("I83").Value = 1498.5
$ws.Range("J83").Value = 997.6667
$ws.Range("K83").Value = 7492.5
$ws.Range("L83").Value = 4988.3335
$ws.Range("M83").Value = -2500.5
$ws.Range("N83").Value = -14972.3335
$ws.Range("H97").Value = 1548.8182
$ws.Range("I97").Value = 1089.2
$ws.Range("J97").Value = 1931.8334
$ws.Range("K97").Value = 1089.2
$ws.Range("L97").Value = 1931.8334
$ws.Range("M97").Value = -593.2
$ws.Range("N97").Value = -2923.8334
$ws.Range("H124").Value = 65000
$ws.Range("J124").Value = 65000
$ws.Range("L124").Value = 65000
$ws.Range("N124").Value = -74820
$ws.Range("H130").Value = 20000
$ws.Range("J130").Value = 20000
$ws.Range("L130").Value = 20000
$ws.Range("N130").Value = -30040

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1181.5
$ws.Range("I22").Value = 455.6
$ws.Range("J22").Value = 1511.4546
$ws.Range("K22").Value = 455.6
$ws.Range("L22").Value = 1511.4546
$ws.Range("M22").Value = -160.6
$ws.Range("N22").Value = -2101.4546
$ws.Range("H27").Value = 1181.5
$ws.Range("I27").Value = 455.6
$ws.Range("J27").Value = 1511.4546
$ws.Range("K27").Value = 455.6
$ws.Range("L27").Value = 1511.4546
$ws.Range("M27").Value = -348.6
$ws.Range("N27").Value = -1725.4546
$ws.Range("H46").Value = 1460.0714
$ws.Range("I46").Value = 960.3333
$ws.Range("K46").Value = 960.3333
$ws.Range("M46").Value = -772.3333
$ws.Range("H99").Value = 40869.1
$ws.Range("I99").Value = 40869.1
$ws.Range("K99").Value = 40869.1
$ws.Range("M99").Value = -37874.1

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 40359.75
$ws.Range("I61").Value = 40359.75
$ws.Range("K61").Value = 40359.75
$ws.Range("M61").Value = -40067.75
$ws.Range("H107").Value = 2298.6365
$ws.Range("I107").Value = 1548.3334
$ws.Range("J107").Value = 3199
$ws.Range("K107").Value = 4645.0002
$ws.Range("L107").Value = 9597
$ws.Range("M107").Value = -2725.0002
$ws.Range("N107").Value = -13437
